$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (hunk 0)
$ws.Cells.Item(15, 8).Value = 1616.6333  # H15: 1591.1476 -> 1616.6333
$ws.Cells.Item(15, 9).Value = 1616.6333  # I15: 1591.1476 -> 1616.6333
$ws.Cells.Item(15, 11).Value = 4849.8999  # K15: 4773.4428 -> 4849.8999
$ws.Cells.Item(15, 13).Value = -4680.8999  # M15: -4604.4428 -> -4680.8999

# Row 33 (hunk 1)
$ws.Cells.Item(33, 8).Value = 412.16666  # H33: 380.33334 -> 412.16666
$ws.Cells.Item(33, 9).Value = 412.16666  # I33: 334.125 -> 412.16666
$ws.Cells.Item(33, 10).Value = 0  # J33: 750 -> 0
$ws.Cells.Item(33, 11).Value = 412.16666  # K33: 334.125 -> 412.16666
$ws.Cells.Item(33, 12).Value = 0  # L33: 750 -> 0
$ws.Cells.Item(33, 13).Value = -183.16666  # M33: -105.125 -> -183.16666
$ws.Cells.Item(33, 14).ClearContents()  # N33: -1208 -> (removed)

# Row 64 (hunk 2)
$ws.Cells.Item(64, 8).Value = 3716.6667  # H64: 3718.182 -> 3716.6667
$ws.Cells.Item(64, 10).Value = 3716.6667  # J64: 3718.182 -> 3716.6667
$ws.Cells.Item(64, 12).Value = 3716.6667  # L64: 3718.182 -> 3716.6667
$ws.Cells.Item(64, 14).Value = -4212.6667  # N64: -4214.182 -> -4212.6667

# Row 67 (hunk 3)
$ws.Cells.Item(67, 8).Value = 3716.6667  # H67: 3718.182 -> 3716.6667
$ws.Cells.Item(67, 10).Value = 3716.6667  # J67: 3718.182 -> 3716.6667
$ws.Cells.Item(67, 12).Value = 3716.6667  # L67: 3718.182 -> 3716.6667
$ws.Cells.Item(67, 14).Value = -5432.6667  # N67: -5434.182 -> -5432.6667

# Row 70 (hunk 4)
$ws.Cells.Item(70, 8).Value = 1400  # H70: 1366.5 -> 1400
$ws.Cells.Item(70, 9).Value = 1100  # I70: 1049.75 -> 1100
$ws.Cells.Item(70, 11).Value = 3300  # K70: 3149.25 -> 3300
$ws.Cells.Item(70, 13).Value = -3030  # M70: -2879.25 -> -3030

# Row 73 (hunk 5)
$ws.Cells.Item(73, 8).Value = 1400  # H73: 1366.5 -> 1400
$ws.Cells.Item(73, 9).Value = 1100  # I73: 1049.75 -> 1100
$ws.Cells.Item(73, 11).Value = 3300  # K73: 3149.25 -> 3300
$ws.Cells.Item(73, 13).Value = -2364  # M73: -2213.25 -> -2364

# Row 74 (hunk 6)
$ws.Cells.Item(74, 8).Value = 3129  # H74: 3277.4 -> 3129
$ws.Cells.Item(74, 9).Value = 2754.8  # I74: 2846.75 -> 2754.8
$ws.Cells.Item(74, 11).Value = 2754.8  # K74: 2846.75 -> 2754.8
$ws.Cells.Item(74, 13).Value = -1818.8  # M74: -1910.75 -> -1818.8

# Row 77 (hunk 7)
$ws.Cells.Item(77, 8).Value = 3129  # H77: 3277.4 -> 3129
$ws.Cells.Item(77, 9).Value = 2754.8  # I77: 2846.75 -> 2754.8
$ws.Cells.Item(77, 11).Value = 13774  # K77: 14233.75 -> 13774
$ws.Cells.Item(77, 13).Value = -9094  # M77: -9553.75 -> -9094

# Row 86 (hunk 8)
$ws.Cells.Item(86, 8).Value = 50  # H86: 0 -> 50
$ws.Cells.Item(86, 9).Value = 50  # I86: 0 -> 50
$ws.Cells.Item(86, 11).Value = 50  # K86: 0 -> 50
$ws.Cells.Item(86, 13).Value = 1073  # M86: None -> 1073

# Row 89 (hunk 9)
$ws.Cells.Item(89, 8).Value = 50  # H89: 0 -> 50
$ws.Cells.Item(89, 9).Value = 50  # I89: 0 -> 50
$ws.Cells.Item(89, 11).Value = 250  # K89: 0 -> 250
$ws.Cells.Item(89, 13).Value = 5366  # M89: None -> 5366

# Row 121 (hunk 10)
$ws.Cells.Item(121, 8).Value = 2043.5416  # H121: 2102.9333 -> 2043.5416
$ws.Cells.Item(121, 10).Value = 2043.5416  # J121: 2102.9333 -> 2043.5416
$ws.Cells.Item(121, 12).Value = 6130.6248  # L121: 6308.7999 -> 6130.6248
$ws.Cells.Item(121, 14).Value = -9624.6248  # N121: -9802.7999 -> -9624.6248

# Row 127 (hunk 11)
$ws.Cells.Item(127, 8).Value = 1057.7142  # H127: 1048.8334 -> 1057.7142
$ws.Cells.Item(127, 9).Value = 1057.7142  # I127: 1048.8334 -> 1057.7142
$ws.Cells.Item(127, 11).Value = 3173.1426  # K127: 3146.5002 -> 3173.1426
$ws.Cells.Item(127, 13).Value = 1786.8574  # M127: 1813.4998 -> 1786.8574

# Row 129 (hunk 12)
$ws.Cells.Item(129, 8).Value = 2863  # H129: 2863.3 -> 2863
$ws.Cells.Item(129, 9).Value = 2734  # I129: 2670.25 -> 2734
$ws.Cells.Item(129, 11).Value = 8202  # K129: 8010.75 -> 8202
$ws.Cells.Item(129, 13).Value = -3202  # M129: -3010.75 -> -3202

# Row 137 (hunk 13)
$ws.Cells.Item(137, 8).Value = 2313.1458  # H137: 2349.0435 -> 2313.1458
$ws.Cells.Item(137, 9).Value = 1626.7333  # I137: 1636.6786 -> 1626.7333
$ws.Cells.Item(137, 11).Value = 4880.199900000001  # K137: 4910.0358 -> 4880.199900000001
$ws.Cells.Item(137, 13).Value = -2330.199900000001  # M137: -2360.0358 -> -2330.199900000001

# Row 138 (hunk 14)
$ws.Cells.Item(138, 8).Value = 5129.077  # H138: 5136.074 -> 5129.077
$ws.Cells.Item(138, 10).Value = 5182.125  # J138: 5187.56 -> 5182.125
$ws.Cells.Item(138, 12).Value = 15546.375  # L138: 15562.68 -> 15546.375
$ws.Cells.Item(138, 14).Value = -25826.375  # N138: -25842.68 -> -25826.375

# Row 141 (hunk 15)
$ws.Cells.Item(141, 8).Value = 0  # H141: 9999 -> 0
$ws.Cells.Item(141, 9).Value = 0  # I141: 9999 -> 0
$ws.Cells.Item(141, 11).Value = 0  # K141: 29997 -> 0
$ws.Cells.Item(141, 13).ClearContents()  # M141: -24817 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 16)
$ws.Cells.Item(32, 8).Value = 21999.4  # H32: 11552.053 -> 21999.4
$ws.Cells.Item(32, 9).Value = 21999.4  # I32: 11552.053 -> 21999.4
$ws.Cells.Item(32, 11).Value = 21999.4  # K32: 11552.053 -> 21999.4
$ws.Cells.Item(32, 13).Value = -21712.4  # M32: -11265.053 -> -21712.4

# Row 97 (hunk 17)
$ws.Cells.Item(97, 8).Value = 333.375  # H97: 252.11111 -> 333.375
$ws.Cells.Item(97, 9).Value = 295.2857  # I97: 252.11111 -> 295.2857
$ws.Cells.Item(97, 10).Value = 600  # J97: 0 -> 600
$ws.Cells.Item(97, 11).Value = 295.2857  # K97: 252.11111 -> 295.2857
$ws.Cells.Item(97, 12).Value = 600  # L97: 0 -> 600
$ws.Cells.Item(97, 13).Value = 200.7143  # M97: 243.88889 -> 200.7143
$ws.Cells.Item(97, 14).Value = -1592  # N97: None -> -1592

# Row 132 (hunk 18)
$ws.Cells.Item(132, 8).Value = 0  # H132: 4399.6665 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 3200 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 4999.5 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 9600 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 14998.5 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: -7070 -> (removed)
$ws.Cells.Item(132, 14).ClearContents()  # N132: -20058.5 -> (removed)

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (hunk 19)
$ws.Cells.Item(94, 8).Value = 1357.4  # H94: 2000 -> 1357.4
$ws.Cells.Item(94, 9).Value = 1349.75  # I94: 2000 -> 1349.75
$ws.Cells.Item(94, 10).Value = 1388  # J94: 0 -> 1388
$ws.Cells.Item(94, 11).Value = 1349.75  # K94: 2000 -> 1349.75
$ws.Cells.Item(94, 12).Value = 1388  # L94: 0 -> 1388
$ws.Cells.Item(94, 13).Value = -898.75  # M94: -1549 -> -898.75
$ws.Cells.Item(94, 14).Value = -2290  # N94: None -> -2290

# Row 134 (hunk 20)
$ws.Cells.Item(134, 8).Value = 4975  # H134: 5148 -> 4975
$ws.Cells.Item(134, 9).Value = 4975  # I134: 4977.6 -> 4975
$ws.Cells.Item(134, 10).Value = 0  # J134: 6000 -> 0
$ws.Cells.Item(134, 11).Value = 14925  # K134: 14932.8 -> 14925
$ws.Cells.Item(134, 12).Value = 0  # L134: 18000 -> 0
$ws.Cells.Item(134, 13).Value = -12390  # M134: -12397.8 -> -12390
$ws.Cells.Item(134, 14).ClearContents()  # N134: -23070 -> (removed)

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 21)
$ws.Cells.Item(31, 8).Value = 2377.84  # H31: 2371.7778 -> 2377.84
$ws.Cells.Item(31, 9).Value = 2377.84  # I31: 2386.5 -> 2377.84
$ws.Cells.Item(31, 10).Value = 0  # J31: 1989 -> 0
$ws.Cells.Item(31, 11).Value = 2377.84  # K31: 2386.5 -> 2377.84
$ws.Cells.Item(31, 12).Value = 0  # L31: 1989 -> 0
$ws.Cells.Item(31, 13).Value = -2082.84  # M31: -2091.5 -> -2082.84
$ws.Cells.Item(31, 14).ClearContents()  # N31: -2579 -> (removed)

# Row 34 (hunk 22)
$ws.Cells.Item(34, 8).Value = 2377.84  # H34: 2371.7778 -> 2377.84
$ws.Cells.Item(34, 9).Value = 2377.84  # I34: 2386.5 -> 2377.84
$ws.Cells.Item(34, 10).Value = 0  # J34: 1989 -> 0
$ws.Cells.Item(34, 11).Value = 2377.84  # K34: 2386.5 -> 2377.84
$ws.Cells.Item(34, 12).Value = 0  # L34: 1989 -> 0
$ws.Cells.Item(34, 13).Value = -2175.84  # M34: -2184.5 -> -2175.84
$ws.Cells.Item(34, 14).ClearContents()  # N34: -2393 -> (removed)

# Row 52 (hunk 23)
$ws.Cells.Item(52, 8).Value = 134825  # H52: 133300 -> 134825
$ws.Cells.Item(52, 10).Value = 139766.67  # J52: 139950 -> 139766.67
$ws.Cells.Item(52, 12).Value = 139766.67  # L52: 139950 -> 139766.67
$ws.Cells.Item(52, 14).Value = -140354.67  # N52: -140538 -> -140354.67

# Row 62 (hunk 24)
$ws.Cells.Item(62, 8).Value = 4639.4  # H62: 4800 -> 4639.4
$ws.Cells.Item(62, 10).Value = 4398.5  # J62: 0 -> 4398.5
$ws.Cells.Item(62, 12).Value = 4398.5  # L62: 0 -> 4398.5
$ws.Cells.Item(62, 14).Value = -5646.5  # N62: None -> -5646.5

# Row 65 (hunk 25)
$ws.Cells.Item(65, 8).Value = 4639.4  # H65: 4800 -> 4639.4
$ws.Cells.Item(65, 10).Value = 4398.5  # J65: 0 -> 4398.5
$ws.Cells.Item(65, 12).Value = 21992.5  # L65: 0 -> 21992.5
$ws.Cells.Item(65, 14).Value = -28232.5  # N65: None -> -28232.5

# Row 139 (hunk 26)
$ws.Cells.Item(139, 8).Value = 120550  # H139: 125000 -> 120550
$ws.Cells.Item(139, 10).Value = 120550  # J139: 125000 -> 120550
$ws.Cells.Item(139, 12).Value = 120550  # L139: 125000 -> 120550
$ws.Cells.Item(139, 14).Value = -130830  # N139: -135280 -> -130830

$ws = $wb.Worksheets.Item("CUL")
# Row 122 (hunk 27)
$ws.Cells.Item(122, 8).Value = 168582  # H122: 144856 -> 168582
$ws.Cells.Item(122, 10).Value = 252123  # J122: 202198.4 -> 252123
$ws.Cells.Item(122, 12).Value = 2269107  # L122: 1819785.6 -> 2269107
$ws.Cells.Item(122, 14).Value = -2274007  # N122: -1824685.6 -> -2274007

# Row 132 (hunk 28)
$ws.Cells.Item(132, 8).Value = 4176.5713  # H132: 4887.2104 -> 4176.5713
$ws.Cells.Item(132, 10).Value = 4064.4167  # J132: 4891.647 -> 4064.4167
$ws.Cells.Item(132, 12).Value = 36579.7503  # L132: 44024.823 -> 36579.7503
$ws.Cells.Item(132, 14).Value = -41639.7503  # N132: -49084.823 -> -41639.7503

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (hunk 29)
$ws.Cells.Item(2, 8).Value = 186.66667  # H2: 177.88889 -> 186.66667
$ws.Cells.Item(2, 9).Value = 186.66667  # I2: 177.88889 -> 186.66667
$ws.Cells.Item(2, 11).Value = 186.66667  # K2: 177.88889 -> 186.66667
$ws.Cells.Item(2, 13).Value = -73.66667000000001  # M2: -64.88889 -> -73.66667000000001

# Row 80 (hunk 30)
$ws.Cells.Item(80, 8).Value = 7231.6665  # H80: 7648.3335 -> 7231.6665
$ws.Cells.Item(80, 9).Value = 5000  # I80: 7500 -> 5000
$ws.Cells.Item(80, 11).Value = 5000  # K80: 7500 -> 5000
$ws.Cells.Item(80, 13).Value = -4002  # M80: -6502 -> -4002

# Row 83 (hunk 31)
$ws.Cells.Item(83, 8).Value = 7231.6665  # H83: 7648.3335 -> 7231.6665
$ws.Cells.Item(83, 9).Value = 5000  # I83: 7500 -> 5000
$ws.Cells.Item(83, 11).Value = 25000  # K83: 37500 -> 25000
$ws.Cells.Item(83, 13).Value = -20008  # M83: -32508 -> -20008

# Row 132 (hunk 32)
$ws.Cells.Item(132, 8).Value = 1899.4  # H132: 2500 -> 1899.4
$ws.Cells.Item(132, 9).Value = 1899.4  # I132: 2000 -> 1899.4
$ws.Cells.Item(132, 10).Value = 0  # J132: 3000 -> 0
$ws.Cells.Item(132, 11).Value = 5698.200000000001  # K132: 6000 -> 5698.200000000001
$ws.Cells.Item(132, 12).Value = 0  # L132: 9000 -> 0
$ws.Cells.Item(132, 13).Value = -3168.200000000001  # M132: -3470 -> -3168.200000000001
$ws.Cells.Item(132, 14).ClearContents()  # N132: -14060 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 33)
$ws.Cells.Item(22, 8).Value = 7199.6  # H22: 7564 -> 7199.6
$ws.Cells.Item(22, 10).Value = 7249.75  # J22: 7658.1665 -> 7249.75
$ws.Cells.Item(22, 12).Value = 7249.75  # L22: 7658.1665 -> 7249.75
$ws.Cells.Item(22, 14).Value = -7839.75  # N22: -8248.166499999999 -> -7839.75

# Row 27 (hunk 34)
$ws.Cells.Item(27, 8).Value = 7199.6  # H27: 7564 -> 7199.6
$ws.Cells.Item(27, 10).Value = 7249.75  # J27: 7658.1665 -> 7249.75
$ws.Cells.Item(27, 12).Value = 7249.75  # L27: 7658.1665 -> 7249.75
$ws.Cells.Item(27, 14).Value = -7463.75  # N27: -7872.1665 -> -7463.75

# Row 46 (hunk 35)
$ws.Cells.Item(46, 8).Value = 449  # H46: 424.16666 -> 449
$ws.Cells.Item(46, 9).Value = 449  # I46: 424.16666 -> 449
$ws.Cells.Item(46, 11).Value = 449  # K46: 424.16666 -> 449
$ws.Cells.Item(46, 13).Value = -261  # M46: -236.16666 -> -261

# Row 55 (hunk 36)
$ws.Cells.Item(55, 8).Value = 420.16666  # H55: 419.8889 -> 420.16666
$ws.Cells.Item(55, 9).Value = 430.5  # I55: 412.85715 -> 430.5
$ws.Cells.Item(55, 10).Value = 399.5  # J55: 444.5 -> 399.5
$ws.Cells.Item(55, 11).Value = 430.5  # K55: 412.85715 -> 430.5
$ws.Cells.Item(55, 12).Value = 399.5  # L55: 444.5 -> 399.5
$ws.Cells.Item(55, 13).Value = -257.5  # M55: -239.85715 -> -257.5
$ws.Cells.Item(55, 14).Value = -745.5  # N55: -790.5 -> -745.5

# Row 100 (hunk 37)
$ws.Cells.Item(100, 8).Value = 2312  # H100: 2816.3333 -> 2312
$ws.Cells.Item(100, 9).Value = 1555.5  # I100: 0 -> 1555.5
$ws.Cells.Item(100, 11).Value = 1555.5  # K100: 0 -> 1555.5
$ws.Cells.Item(100, 13).Value = -1014.5  # M100: None -> -1014.5

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (hunk 38)
$ws.Cells.Item(126, 8).Value = 2215.5  # H126: 2264.4211 -> 2215.5
$ws.Cells.Item(126, 9).Value = 2253.158  # I126: 2264.4211 -> 2253.158
$ws.Cells.Item(126, 10).Value = 1500  # J126: 0 -> 1500
$ws.Cells.Item(126, 11).Value = 6759.474  # K126: 6793.263300000001 -> 6759.474
$ws.Cells.Item(126, 12).Value = 4500  # L126: 0 -> 4500
$ws.Cells.Item(126, 13).Value = -4289.474  # M126: -4323.263300000001 -> -4289.474
$ws.Cells.Item(126, 14).Value = -9440  # N126: None -> -9440
